$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-6
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45208
}
